$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Add the new tracked-file entry as a new row in the log table (row 49):
#   B49 = commit date, C49 = file path, D49 = lines changed
# Copy the date formatting from an existing date cell in column B so the
# new cell reuses the same number format style as the rest of the table.
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B49").Value = 45916
$ws.Range("C49").Value = "Tokyo-Garden-Restaurant/Dokumentacja projektu Tokyo Garden.docx"
$ws.Range("D49").Value = 131

# Restore the sheet's scroll position / selection as left by the author.
$ws.Application.Goto($ws.Range("A3"))
$ws.Range("E51").Select()
